$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individualios_rek")

$ws.Range("B2").Value = "Receptas: AVOKADŲ IR PESTO MAKARONŲ SALOTOS"
$ws.Range("B3").Value = "Produktas: Baklažanai"
$ws.Range("B4").Value = "Receptas: KREMINIAI GRYBŲ MAKARONAI"
$ws.Range("B5").Value = "Produktas: Rabarbarai"
